$d = $word.ActiveDocument

# Locate the "Link GitHub" heading paragraph, the empty paragraph right
# after it, and the "Exercício" heading paragraph that follows.
$linkParaIndex = 0
$emptyParaIndex = 0
$exercicioParaIndex = 0

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = [string]$p.Range.Text

    if ($text -match "Link GitHub") {
        $linkParaIndex = $i
    }
    elseif ($linkParaIndex -gt 0 -and $emptyParaIndex -eq 0 -and $text.Trim() -eq "") {
        $emptyParaIndex = $i
    }
    elseif ($emptyParaIndex -gt 0 -and $exercicioParaIndex -eq 0 -and $text -match "Exerc") {
        $exercicioParaIndex = $i
        break
    }
}

# Remove the stray empty paragraph (<w:p/>) that sits between the two
# headings.
$emptyPara = $d.Paragraphs.Item($emptyParaIndex)
$emptyPara.Range.Delete()

# The "Exercício" paragraph has shifted up by one after the delete.
$targetIndex = $exercicioParaIndex - 1
$targetPara = $d.Paragraphs.Item($targetIndex)

$paraStart = $targetPara.Range.Start
$paraEnd = $targetPara.Range.End

# Clear the "Exercício" text but keep the paragraph mark itself.
$textRange = $d.Range($paraStart, $paraEnd - 1)
$textRange.Text = ""

# Drop the Heading ("Ttulo1") style so the paragraph falls back to Normal,
# matching the target markup (no <w:pPr> at all).
$targetPara = $d.Paragraphs.Item($targetIndex)
$targetPara.Style = $d.Styles.Item("Normal")

# Insert the GitHub repository link as a real hyperlink whose display text
# is the URL itself.
$url = "https://github.com/faculdade-infnet/IV-2-C_sharp/tree/main/TP2"
$insertRange = $d.Range($paraStart, $paraStart)
$d.Hyperlinks.Add($insertRange, $url, [Type]::Missing, [Type]::Missing, $url) | Out-Null
